$d = $word.ActiveDocument

# Locate the paragraph whose text is the "Artigos de revistas..." bibliography
# entry; the three paragraphs right after it (a blank paragraph, the
# "Ver no Jupiter..." line and the "© 2020 ..." footer line) are being
# removed by this edit, while the paragraph that follows them survives.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Artigos de revistas especializadas*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -ge 1) {
    $firstParaToRemove = $d.Paragraphs.Item($anchorIndex + 1)
    $lastParaToRemove = $d.Paragraphs.Item($anchorIndex + 3)

    $deleteRange = $d.Range($firstParaToRemove.Range.Start, $lastParaToRemove.Range.End)
    $deleteRange.Delete()
}
